# Applies the cryptos.xlsx data refresh described in the commit
# "Updated cryptos list on Sun May 26 11:44:48 UTC 2024 with GitHub Actions".
# Column D (Price) holds plain text (values such as "69.097.64" use dots as
# thousands separators, not decimal points), so each Price cell is forced to
# the Text number format before the value is written - this keeps Excel from
# auto-converting it to a number and silently dropping precision / trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.097.64"
$ws.Range("E2").Value = "  -0.49%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.807.79"
$ws.Range("E3").Value = "  +1.45%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.56"
$ws.Range("E5").Value = "  -0.48%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.91"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.805.59"
$ws.Range("E7").Value = "  +1.48%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("E9").Value = "  -0.59%  "

# Row 10
$ws.Range("E10").Value = "  +1.70%  "

# Row 11
$ws.Range("E11").Value = "  -1.28%  "

# Row 12
$ws.Range("E12").Value = "  -0.08%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.26"
$ws.Range("E13").Value = "  -3.17%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000246"
$ws.Range("E14").Value = "  -1.32%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.447.88"
$ws.Range("E15").Value = "  +1.46%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.806.83"
$ws.Range("E16").Value = "  +1.24%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.219.66"
$ws.Range("E17").Value = "  -0.28%  "

# Row 18
$ws.Range("E18").Value = "  +1.87%  "

# Row 19
$ws.Range("E19").Value = "  -0.24%  "

# Row 20
$ws.Range("E20").Value = "  +1.29%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.41"
$ws.Range("E21").Value = "  +4.43%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "488.13"
$ws.Range("E22").Value = "  -1.57%  "

# Row 23
$ws.Range("E23").Value = "  -0.72%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000158"
$ws.Range("E24").Value = "  +3.57%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.84"
$ws.Range("E25").Value = "  -0.75%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.26"
$ws.Range("E26").Value = "  -3.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.23"
$ws.Range("E27").Value = "  -1.22%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("E28").Value = "  -2.40%  "

# Row 29
$ws.Range("E29").Value = "  +0.03%  "

# Row 30
$ws.Range("E30").Value = "  -1.15%  "

# Row 31
$ws.Range("E31").Value = "  -0.27%  "

# Row 32
$ws.Range("E32").Value = "  -4.14%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.958.19"
$ws.Range("E33").Value = "  +1.43%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.81"
$ws.Range("E34").Value = "  -0.31%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.753.37"
$ws.Range("E35").Value = "  +1.81%  "

# Row 36
$ws.Range("E36").Value = "  -1.53%  "

# Row 37
$ws.Range("E37").Value = "  +5.46%  "

# Row 38
$ws.Range("E38").Value = "  +0.37%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.91"
$ws.Range("E39").Value = "  +0.20%  "

# Row 40
$ws.Range("E40").Value = "  -0.04%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.320"
$ws.Range("E41").Value = "  -1.65%  "

# Row 42
$ws.Range("E42").Value = "  -0.57%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.60"
$ws.Range("E43").Value = "  -0.19%  "

# Row 44
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.99"
$ws.Range("E44").Value = "  -0.29%  "

# Row 45
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "425.82"
$ws.Range("E45").Value = "  -2.79%  "

# Row 46
$ws.Range("E46").Value = "  +0.00%  "

# Row 47
$ws.Range("E47").Value = "  -1.20%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.835.60"
$ws.Range("E48").Value = "  +1.04%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.34"
$ws.Range("E49").Value = "  +0.36%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "39.58"
$ws.Range("E50").Value = "  -2.62%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0350"
$ws.Range("E51").Value = "  -1.25%  "
